# Append a new data row (row 48) to each of the 4 worksheets, mirroring the
# structure/format of the existing rows (e.g. row 47) but with an updated
# timestamp in column A.

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"
$newTime = 45834.43858796296

# --- Sheet 1: DE_LFT_#1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A48").Value = $newTime
$ws.Range("A48").NumberFormat = $dateFormat
$ws.Range("B48").Value = "0x01,0x7c"
$ws.Range("C48").Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Range("D48").Value = "0x01,0x64"
$ws.Range("E48").Value = "0x14"
$ws.Range("F48").Value = 380
$ws.Range("G48").Value = [double]"7.598631275147109e+23"
$ws.Range("H48").Value = 356
$ws.Range("I48").Value = 14

# --- Sheet 2: DE_LFT_#2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A48").Value = $newTime
$ws.Range("A48").NumberFormat = $dateFormat
$ws.Range("B48").Value = "0x01,0x7c"
$ws.Range("C48").Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Range("D48").Value = "0x01,0x64"
$ws.Range("E48").Value = "0xe"
$ws.Range("F48").Value = 380
$ws.Range("G48").Value = [double]"5.68432987514711e+23"
$ws.Range("H48").Value = 356
$ws.Range("I48").Value = 14

# --- Sheet 3: DE_PLT_#1 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A48").Value = $newTime
$ws.Range("A48").NumberFormat = $dateFormat
$ws.Range("B48").Value = "0x00,0x82"
$ws.Range("C48").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Range("D48").Value = "0x00,0x7F"
$ws.Range("E48").Value = "0x7"
$ws.Range("F48").Value = 130
$ws.Range("G48").Value = [double]"5.68631262647114e+23"
$ws.Range("H48").Value = 127
$ws.Range("I48").Value = 7

# --- Sheet 4: DE_PLT_#2 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("A48").Value = $newTime
$ws.Range("A48").NumberFormat = $dateFormat
$ws.Range("B48").Value = "0x00,0x82"
$ws.Range("C48").Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Range("D48").Value = "0x00,0x7E"
$ws.Range("E48").Value = "0x3"
$ws.Range("F48").Value = 130
$ws.Range("G48").Value = [double]"9.85046333984776e+23"
$ws.Range("H48").Value = 126
$ws.Range("I48").Value = 3

Write-Output "Row 48 appended to all 4 sheets"
